$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'311.31"
$ws.Range("E2").Value = "'1.07%"
$ws.Range("D3").Value = "'39.41"
$ws.Range("E3").Value = "'2.03%"
$ws.Range("D4").Value = "'5.140"
$ws.Range("E4").Value = "'0.87%"
$ws.Range("D5").Value = "'0.08147"
$ws.Range("E5").Value = "'0.15%"
$ws.Range("D6").Value = "'1.993"
$ws.Range("E6").Value = "'1.28%"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'8.137"
$ws.Range("E7").Value = "'2.58%"
$ws.Range("B8").Value = "GateToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D8").Value = "'4.234"
$ws.Range("E8").Value = "'0.76%"
$ws.Range("D9").Value = "'0.9265"
$ws.Range("E9").Value = "'-0.18%"
$ws.Range("D10").Value = "'0.1405"
$ws.Range("E10").Value = "'-2.61%"
$ws.Range("D11").Value = "'0.1929"
$ws.Range("E11").Value = "'-1.51%"
$ws.Range("D12").Value = "'0.09034"
$ws.Range("E12").Value = "'-0.71%"
$ws.Range("D13").Value = "'0.03510"
$ws.Range("E13").Value = "'0.00%"
$ws.Range("D14").Value = "'0.09816"
$ws.Range("E14").Value = "'-0.07%"
$ws.Range("D15").Value = "'0.001398"
$ws.Range("E15").Value = "'-0.49%"
$ws.Range("D16").Value = "'0.006023"
$ws.Range("E16").Value = "'-0.93%"
$ws.Range("E17").Value = "'1.41%"
$ws.Range("D18").Value = "'3.360"
$ws.Range("E18").Value = "'-2.90%"
$ws.Range("D19").Value = "'0.3456"
$ws.Range("E19").Value = "'0.31%"
$ws.Range("D20").Value = "'0.1311"
$ws.Range("E20").Value = "'-1.73%"
$ws.Range("D21").Value = "'4.642"
$ws.Range("E21").Value = "'-3.77%"
$ws.Range("D22").Value = "'0.2422"
$ws.Range("E22").Value = "'0.80%"
$ws.Range("D23").Value = "'0.04366"
$ws.Range("E23").Value = "'-1.53%"
$ws.Range("E24").Value = "'0.89%"
$ws.Range("D25").Value = "'0.004883"
$ws.Range("E25").Value = "'0.72%"
$ws.Range("D26").Value = "'0.0001299"
$ws.Range("E26").Value = "'-0.22%"
$ws.Range("D27").Value = "'0.0003997"
$ws.Range("E27").Value = "'-10.13%"
$ws.Range("D39").Value = "'0.02141"
$ws.Range("E39").Value = "'1.93%"
$ws.Range("D40").Value = "'0.05186"
$ws.Range("E40").Value = "'1.26%"
$ws.Range("D41").Value = "'0.007426"
$ws.Range("E41").Value = "'-0.83%"
$ws.Range("D42").Value = "'0.009864"
$ws.Range("E42").Value = "'-2.82%"
$ws.Range("D43").Value = "'0.1369"
$ws.Range("E43").Value = "'0.46%"
$ws.Range("D44").Value = "'0.002119"
$ws.Range("E44").Value = "'-1.13%"
$ws.Range("D45").Value = "'0.008995"
$ws.Range("E45").Value = "'-4.59%"
$ws.Range("D46").Value = "'0.00006409"
$ws.Range("E46").Value = "'2.90%"
$ws.Range("E47").Value = "'-0.20%"
$ws.Range("D48").Value = "'0.0009992"
$ws.Range("E48").Value = "'-37.62%"
$ws.Range("D49").Value = "'0.002542"
$ws.Range("E49").Value = "'-16.94%"
$ws.Range("D50").Value = "'0.00002099"
$ws.Range("E50").Value = "'-0.20%"
$ws.Range("D51").Value = "'0.0001999"
$ws.Range("E51").Value = "'-0.20%"
